$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- Row 17: begin group / cholera_cases / Cholera Cases ---
$ws.Range("A17").Value = "begin group"
$ws.Range("B17").Value = "cholera_cases"
$ws.Range("C17").Value = "Cholera Cases"

# --- Row 18: integer / num_screened / Number Of Persons Screened For Cholera / yes // numbers ---
$ws.Range("A18").Value = "integer"
$ws.Range("B18").Value = "num_screened"
$ws.Range("C18").Value = "Number Of Persons Screened For Cholera"
$ws.Range("D18").Value = "yes"
$ws.Range("F18").Value = "numbers"

# --- Row 19: integer / num_referred / Number Of Presumptive Cholera Persons Referred For Diagnosis / yes // numbers ---
$ws.Range("A19").Value = "integer"
$ws.Range("B19").Value = "num_referred"
$ws.Range("C19").Value = "Number Of Presumptive Cholera Persons Referred For Diagnosis"
$ws.Range("D19").Value = "yes"
$ws.Range("F19").Value = "numbers"

# --- Row 20: integer / num_referred_reached / Number Of Referred Persons Who Reached Health Facility / yes / ${num_referred} >0 / numbers ---
$ws.Range("A20").Value = "integer"
$ws.Range("B20").Value = "num_referred_reached"
$ws.Range("C20").Value = "Number Of Referred Persons Who Reached Health Facility"
$ws.Range("D20").Value = "yes"
$ws.Range("E20").Value = '${num_referred} >0'
$ws.Range("F20").Value = "numbers"

# --- Row 21: integer / num_confirmed_cases / Number Of Confirmed Cholera Cases At Health Facility / yes // numbers ---
$ws.Range("A21").Value = "integer"
$ws.Range("B21").Value = "num_confirmed_cases"
$ws.Range("C21").Value = "Number Of Confirmed Cholera Cases At Health Facility"
$ws.Range("D21").Value = "yes"
$ws.Range("F21").Value = "numbers"

# --- Row 22: integer / num_deaths / Number Of Deaths Due To Cholera In The Month / yes // numbers ---
$ws.Range("A22").Value = "integer"
$ws.Range("B22").Value = "num_deaths"
$ws.Range("C22").Value = "Number Of Deaths Due To Cholera In The Month"
$ws.Range("D22").Value = "yes"
$ws.Range("F22").Value = "numbers"

# --- Row 23: end group ---
$ws.Range("A23").Value = "end group"

# --- Placeholder empty-but-styled cells that round out row 17 (mirrors D7/F7 elsewhere in the sheet) ---
$ws.Range("D17").NumberFormat = "General"
$ws.Range("F17").NumberFormat = "General"

# --- Apply the sheet's normal style ("1") to all newly written cells, column by column ---
$ws.Range("A1").Copy()
$ws.Range("A17:A23").PasteSpecial(-4122)
$ws.Range("B17:B22").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C19:C22").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("D19:D22").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("F19:F22").PasteSpecial(-4122)

# --- Row 18's C/D/F use the shaded/highlighted style ("2") seen elsewhere in the form (e.g. C15) ---
$ws.Range("C15").Copy()
$ws.Range("C18:D18").PasteSpecial(-4122)
$ws.Range("F18").PasteSpecial(-4122)

$excel.CutCopyMode = $false
